# Weekly update: a new price record is added at the top of the "Poroto
# granado" (Comercializadora del Agro de Limarí) series, and another new
# record is inserted further down in the chronological list. All of the
# existing rows shift down to make room (Excel "insert entire row"
# semantics), which is why most of the sheet's D/J/K/L/M/N/P/Q values
# appear to move down by one (or two) rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the first new record at row 66 (top of the data block) -------
$ws.Rows.Item(66).Insert()

$ws.Range("A66").Value = 2
$ws.Range("B66").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C66").Value = "Coquimbo"
$ws.Range("D66").Value = 45280
$ws.Range("E66").Value = 4
$ws.Range("F66").Value = 100112030
$ws.Range("G66").Value = "Poroto granado"
$ws.Range("H66").Value = "Sin especificar"
$ws.Range("I66").Value = "Primera"
$ws.Range("J66").Value = 300
$ws.Range("K66").Value = 25000
$ws.Range("L66").Value = 26000
$ws.Range("M66").Value = 25500
$ws.Range("N66").Value = "`$/caja 15 kilos"
$ws.Range("O66").Value = "Provincia de Limarí"
$ws.Range("P66").Value = 1700
$ws.Range("Q66").Value = 15
$ws.Range("R66").Value = "Hortaliza"

# --- Insert the second new record at (new) row 109 -----------------------
$ws.Rows.Item(109).Insert()

$ws.Range("A109").Value = 2
$ws.Range("B109").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C109").Value = "Coquimbo"
$ws.Range("D109").Value = 44902
$ws.Range("E109").Value = 4
$ws.Range("F109").Value = 100112030
$ws.Range("G109").Value = "Poroto granado"
$ws.Range("H109").Value = "Sin especificar"
$ws.Range("I109").Value = "Primera"
$ws.Range("J109").Value = 500
$ws.Range("K109").Value = 50000
$ws.Range("L109").Value = 52000
$ws.Range("M109").Value = 51000
$ws.Range("N109").Value = "`$/malla 25 kilos"
$ws.Range("O109").Value = "Provincia de Limarí"
$ws.Range("P109").Value = 2040
$ws.Range("Q109").Value = 25
$ws.Range("R109").Value = "Hortaliza"
